# Append the new resale-numbers row (row 36) reported at 2025-01-21 22:36:28.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Date, Time, Weekday, Week) -----------------------------
# A plain Value assignment of a date-/number-looking string (e.g.
# "2025-01-21" or "03") gets auto-coerced into a real date/number by Excel.
# To land these as literal text (matching the rest of the column, which has
# no explicit number format) we build each value via a text formula and then
# convert it to a static value with Copy / PasteSpecial (values only) -- this
# avoids leaving any NumberFormat/style residue behind on the cells.
$ws.Range("A36").Formula = "=""2025-01-21"""
$ws.Range("B36").Formula = "=""22:36:28"""
$ws.Range("C36").Formula = "=""Tuesday"""
$ws.Range("D36").Formula = "=""03"""

$ws.Range("A36:D36").Copy()
$ws.Range("A36:D36").PasteSpecial(-4163)  # xlPasteValues

# --- Numeric columns (per-city resale counts, -1 = no data) ---------------
$ws.Cells.Item(36, 5).Value  = 126320   # Beijing
$ws.Cells.Item(36, 6).Value  = 142142   # Guangzhou
$ws.Cells.Item(36, 7).Value  = 168671   # Suzhou
$ws.Cells.Item(36, 8).Value  = 158605   # Hangzhou
$ws.Cells.Item(36, 9).Value  = -1       # Nanjing
$ws.Cells.Item(36, 10).Value = 143006   # Xi_an
$ws.Cells.Item(36, 11).Value = -1       # Chengdu
$ws.Cells.Item(36, 12).Value = -1       # Chongqing
$ws.Cells.Item(36, 13).Value = 192267   # Tianjin
$ws.Cells.Item(36, 14).Value = 115689   # Hefei
$ws.Cells.Item(36, 15).Value = 45618    # Fuzhou
$ws.Cells.Item(36, 16).Value = 28471    # Xiamen
$ws.Cells.Item(36, 17).Value = 65812    # Changsha
$ws.Cells.Item(36, 18).Value = -1       # Shanghai
$ws.Cells.Item(36, 19).Value = 48745    # Shenzhen
$ws.Cells.Item(36, 20).Value = -1       # Wuhan
